$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet tab from "CVMod8 V2" to "BOM"
$ws.Name = "BOM"

# Fix the bi-color LED package naming: "LED 3mm Flat Bicolor" -> "LED 3mm Dome Bicolor"
$ws.Range("C14").Value = "LED 3mm Dome Bicolor"

# Update the active selection to C14 (matches the saved selection state in the file)
$ws.Range("C14").Select()
